# This workbook is a weekly price-log sheet for "Haba" (fava bean) at the
# "Macroferia Regional de Talca" market. Each data row is one weekly
# observation; rows are ordered by date. A new weekly observation needs to
# be inserted as the new row 27 (right after row 26, 2020-11-02 / 44159),
# and every existing observation from the old row 27 onward shifts down by
# one row (old row 27 -> new row 28, ..., old row 56 -> new row 57).
#
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across every data row in this
# sheet (same market/category/unit/etc.), so after the shift we just need
# to copy those constant columns into the freshly inserted row, plus carry
# the "Origen" (O) value forward from what used to be row 27 (it doesn't
# change), and write the brand-new D/J/K/L/M/P values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "Origen" (column O = 15) text that used to live in row 27 -
# it stays attached to the same logical record after the insert shifts it
# down to row 28, and the new row 27 keeps that same origin value too.
$origen27 = $ws.Cells.Item(27, 15).Value2

# Insert a new blank row above the current row 27; everything from the old
# row 27 downward (through row 56) shifts down to rows 28-57, and the
# sheet's used range / dimension grows from R56 to R57 automatically.
$ws.Rows.Item(27).Insert()

# The inserted row is blank except for some inherited formatting. Copy the
# columns that are constant for every record in this data set from the
# row right below (the shifted former row 27, now row 28).
$constantCols = 1,2,3,5,6,7,8,9,14,17,18
foreach ($col in $constantCols) {
    $ws.Cells.Item(27, $col).Value = $ws.Cells.Item(28, $col).Value2
}

# New weekly observation values for the inserted row 27.
$ws.Cells.Item(27, 4).Value = 44512     # Fecha
$ws.Cells.Item(27, 10).Value = 600      # Volumen
$ws.Cells.Item(27, 11).Value = 6000     # Precio minimo
$ws.Cells.Item(27, 12).Value = 6000     # Precio maximo
$ws.Cells.Item(27, 13).Value = 6000     # Precio promedio ponderado
$ws.Cells.Item(27, 15).Value = $origen27 # Origen (unchanged)
$ws.Cells.Item(27, 16).Value = 240      # Precio $/Kg
